$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85, shifting existing rows 85:185 down to 86:186.
$ws.Rows.Item(85).Insert()

# Populate the newly inserted row 85 with this week's new data point.
$ws.Range("A85").Value = 11
$ws.Range("B85").Value = "Vega Monumental Concepción"
$ws.Range("C85").Value = "Bíobío"
$ws.Range("D85").Value = 45167
$ws.Range("E85").Value = 8
$ws.Range("F85").Value = 100112001
$ws.Range("G85").Value = "Berenjena"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 80
$ws.Range("K85").Value = 10000
$ws.Range("L85").Value = 10000
$ws.Range("M85").Value = 10000
$ws.Range("N85").Value = "$/caja 50 unidades"
$ws.Range("O85").Value = "Región de Arica y Parinacota"
$ws.Range("P85").Value = 200
$ws.Range("Q85").Value = 50
$ws.Range("R85").Value = "Hortaliza"
